# Update Month/Quarter with 2019, Personal Income
#
# The sheet holds a Month/Quarter lookup table starting in 2020. This
# edit prepends the twelve months of 2019 and also backfills a missing
# March 2020 row (the original data jumped straight from Feb to Apr 2020).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 13 new data rows are needed: Jan 2019 .. Dec 2019, plus the missing
# Mar 2020 row. Insert 13 blank rows right after the header row, pushing
# all existing data down (row 2 -> row 15, etc).
$ws.Rows("2:14").Insert()

# The inserted rows come back with the default/blank style. Copy the
# number-format / alignment of the (now shifted) first original data
# row onto the new block so every new cell matches the rest of the
# table (date format in col A, centered text style in col B).
$ws.Range("A15:B15").Copy()
$ws.Range("A2:B14").PasteSpecial(-4122)

# Fill in the new rows: the 12 months of 2019 plus the backfilled
# March 2020 row.
$dates = @(43466,43497,43525,43556,43586,43617,43647,43678,43709,43739,43770,43800,43891)
$quarters = @("2019Q1","2019Q1","2019Q1","2019Q2","2019Q2","2019Q2","2019Q3","2019Q3","2019Q3","2019Q4","2019Q4","2019Q4","2020Q1")

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $quarters[$i]
}

# The first data row ends up with a slightly taller custom row height.
$ws.Rows(2).RowHeight = 17

# Leave the selection where the author left it.
[void]$ws.Range("C29").Select()
